$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.689.66"
$ws.Range("E2").Value = "  -0.69%  "

$ws.Range("D3").Value = "3.401.00"
$ws.Range("E3").Value = "  -0.66%  "

$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "408.76"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.21%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "127.10"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.00%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.617"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.41%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.999"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.05%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.721"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.29%  "

$ws.Range("E10").Value = "  -9.80%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "42.32"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.00%  "

$ws.Range("B12").Value = "Polkadot"
$ws.Range("C12").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "9.10"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.45%  "

$ws.Range("B13").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C13").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D13").Value = "3.935.12"
$ws.Range("E13").Value = "  -0.76%  "

$ws.Range("E14").Value = "  -0.36%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000208"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -8.60%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "20.29"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.72%  "

$ws.Range("D17").Value = "3.405.61"
$ws.Range("E17").Value = "  -0.60%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.08"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.74%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.23"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.27%  "

$ws.Range("D20").Value = "61.734.35"
$ws.Range("E20").Value = "  -0.46%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "484.56"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +20.52%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "89.43"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.27%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.21"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.62%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.14"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.70%  "

$ws.Range("E25").Value = "  +0.86%  "

$ws.Range("B26").Value = "Filecoin"
$ws.Range("C26").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.48"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +8.44%  "

$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "33.26"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.35%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "4.80"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.06%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.85"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.09%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.65"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.88%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "11.78"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.02%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.168"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.82%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.112"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -6.02%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "40.91"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -5.06%  "

$ws.Range("E35").Value = "  -0.64%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "55.78"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.48%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0484"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.92%  "

$ws.Range("E38").Value = "  +0.11%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.328"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.96%  "

$ws.Range("B40").Value = "Stacks"
$ws.Range("C40").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.97"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.91%  "

$ws.Range("B41").Value = "Monero"
$ws.Range("C41").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "147.60"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.89%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.32"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.58%  "

$ws.Range("E43").Value = "  -0.60%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.06"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.04%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.54"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +5.44%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.17"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.52%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.33"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +16.44%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "16.35"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.97%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "21.94"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.28%  "

$ws.Range("B50").Value = "BitcoinSV"
$ws.Range("C50").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "112.73"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +14.39%  "

$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.142"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +7.87%  "
